# Update the title of slide 31 from "Bài tập 1" to "Bài tập 12.1",
# matching the author's edit which ends up as three runs:
#   "Bài " / "tập " / "12.1"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(31)
$shp = $s.Shapes.Title
$tr = $shp.TextFrame.TextRange

# Original text is "Bài tập 1" (9 chars). Keep the first run ("Bài ",
# chars 1-4) untouched and replace the remainder ("tập 1", chars 5-9)
# with the new wording, split into separate runs the way PowerPoint
# would after the user retyped "tập " and then "12.1".
$tr.Characters(5, 5).Delete() | Out-Null
$tr.InsertAfter("tập ") | Out-Null
$tr.InsertAfter("12.1") | Out-Null
